$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 78, pushing the existing rows 78-82 down to 79-83.
$ws.Rows.Item(78).Insert()

# Copy the (now shifted) formatting/content of row 79 into the new row 78
# for the columns that stay constant across this product block, then set
# the row-specific values that differ (date, volume, prices, price/kg).
$ws.Range("A78").Value = 10
$ws.Range("B78").Value = "Vega Modelo de Temuco"
$ws.Range("C78").Value = "La Araucanía"
$ws.Range("D78").Value = 44753
$ws.Range("D78").NumberFormat = $ws.Range("D79").NumberFormat
$ws.Range("E78").Value = 9
$ws.Range("F78").Value = "Fruta"
$ws.Range("G78").Value = 100108
$ws.Range("H78").Value = "Tropicales y subtropicales"
$ws.Range("I78").Value = 100108007
$ws.Range("J78").Value = "Coco"
$ws.Range("K78").Value = "Sin especificar"
$ws.Range("L78").Value = "Primera"
$ws.Range("M78").Value = 20
$ws.Range("N78").Value = 25000
$ws.Range("O78").Value = 25000
$ws.Range("P78").Value = 25000
$ws.Range("Q78").Value = "$/malla 20 unidades"
$ws.Range("R78").Value = "Perú"
$ws.Range("S78").Value = 1250
$ws.Range("T78").Value = 20
